$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / "Changed") date value from 45742 to 45743
# for all data rows (rows 2 through 45).
for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2() -eq 45742) {
        $cell.Value = 45743
    }
}
